# feat(compare): proof of concept
#
# Duplicate Sheet1 into a new "Sheet3" appended at the end of the workbook,
# select cell F37 on it, and make it the active (visible/selected) sheet.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Sheet1")
$source.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Sheet3"

$newSheet.Activate()
$newSheet.Range("F37").Select()

# Best-effort: reflect the maximized/resized window geometry recorded in the
# target workbook (xWindow/yWindow/windowWidth/windowHeight). Some hosts
# don't persist this cosmetic state, so failures here are non-fatal.
try {
    $win = $excel.ActiveWindow
    $win.WindowState = -4137
    $win.Left = -120
    $win.Top = -120
    $win.Width = 29040
    $win.Height = 15840
} catch {
}
